$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 previously held the shared string "R40"; it now holds the new
# string "1" (appended as a new shared-string entry).
$ws.Range("B11").Value = "1"
